$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Submissions")

$ws.Range("C8").Value = "None"
$ws.Range("D8").Value = "<1 year"
$ws.Range("F8").Value = "Portugal"
$ws.Range("H8").Value = '$0-$500'
$ws.Range("I8").Value = "Rust, DAO Contributor"

# T8 holds a numeric-looking ID that must stay text (matches source which is inlineStr)
$ws.Range("T8").NumberFormat = "@"
$ws.Range("T8").Value = "1697931"

$ws.Range("Y8").Value = "2025-12-19 10:13:32"
$ws.Range("Z8").Value = "Queued"

$wb.Save()
